# Grouping cells takes into account already merged cells + possibility to
# disable grouping of blank cells.
#
# This reproduces the change to TestGroupResultHorizontal.xlsx:
#  - extend two existing horizontal merges by one column (E2:F2 -> E2:G2,
#    B5:C5 -> B5:D5)
#  - add six new data rows (7-12) with grouped/merged header-like cells

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend existing merges ------------------------------------------------
$ws.Range("E2:G2").Merge()
$ws.Range("B5:D5").Merge()

# --- New row 7 ---------------------------------------------------------
$ws.Range("A7").Value = "One"
$ws.Range("B7:I7").Merge()
$ws.Range("B7").Value = "One"
$ws.Range("J7").Value = "One"

# --- New row 8 ---------------------------------------------------------
$ws.Range("A8").Value = "Two"
$ws.Range("B8").Value = "Two"
$ws.Range("C8:I8").Merge()
$ws.Range("C8").Value = "One"
$ws.Range("J8").Value = "One"

# --- New row 9 ---------------------------------------------------------
$ws.Range("A9").Value = "One"
$ws.Range("B9:H9").Merge()
$ws.Range("B9").Value = "One"
$ws.Range("I9").Value = "Two"
$ws.Range("J9").Value = "Two"

# --- New row 10 --------------------------------------------------------
$ws.Range("A10").Value = "One"
$ws.Range("B10:G10").Merge()
$ws.Range("B10").Value = "One"
$ws.Range("H10:I10").Merge()
$ws.Range("H10").Value = "Two"
$ws.Range("J10").Value = "Two"

# --- New row 11 --------------------------------------------------------
$ws.Range("A11").Value = "Two"
$ws.Range("B11").Value = "Two"
$ws.Range("C11:H11").Merge()
$ws.Range("C11").Value = "One"
$ws.Range("I11").Value = "Two"
$ws.Range("J11").Value = "Two"

# --- New row 12 --------------------------------------------------------
$ws.Range("A12:B12").Merge()
$ws.Range("A12").Value = "One"
$ws.Range("C12:G12").Merge()
$ws.Range("C12").Value = "One"
$ws.Range("H12:I12").Merge()
$ws.Range("H12").Value = "Two"
